$wb = $excel.ActiveWorkbook

# Sheet "Manila Philippines": row 4, columns L:W set to 0
$ws = $wb.Worksheets.Item("Manila Philippines")
$ws.Range("L4:W4").Value = 0

# Sheet "Milwaukee Pmc Hq Wisconsin": L7 cleared
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("L7").ClearContents()

# Sheet "Milwaukee Wisconsin": L5 cleared
$ws = $wb.Worksheets.Item("Milwaukee Wisconsin")
$ws.Range("L5").ClearContents()

# Sheet "South Beloit Gardner St Illino": L7 cleared
$ws = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$ws.Range("L7").ClearContents()

# Sheet "Rock Road Radford Virginia": L3 cleared
$ws = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws.Range("L3").ClearContents()
